$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185; this shifts rows 185-239 down to 186-240
# (carrying all their existing values/formatting with them), matching the
# target diff where every record from the old row 185 onward moves down by
# one row and a brand new record appears at row 185.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record's data.
$ws.Range("A185").Value = 10
$ws.Range("B185").Value = "Vega Modelo de Temuco"
$ws.Range("C185").Value = "La Araucanía"
$ws.Range("D185").Value = 44722
$ws.Range("D185").NumberFormat = $ws.Range("D186").NumberFormat
$ws.Range("E185").Value = 9
$ws.Range("F185").Value = "Fruta"
$ws.Range("G185").Value = 100103
$ws.Range("H185").Value = "Frutos de hueso (carozo)"
$ws.Range("I185").Value = 100103002
$ws.Range("J185").Value = "Ciruela"
$ws.Range("K185").Value = "Pink Delight"
$ws.Range("L185").Value = "Primera"
$ws.Range("M185").Value = 55
$ws.Range("N185").Value = 11000
$ws.Range("O185").Value = 12000
$ws.Range("P185").Value = 11455
$ws.Range("Q185").Value = "$/bandeja 18 kilos granel"
$ws.Range("R185").Value = "Región de O'Higgins"
$ws.Range("S185").Value = 636
$ws.Range("T185").Value = 18
